# Applies the codeforiati SectorGroup column re-mapping:
#   new D (category-name) = old E
#   new E (group-name)    = old G
#   new F (category-code) = old F   (unchanged)
#   new G (group-code)    = old D
# across every row of the sheet (including the header row, whose labels
# follow the exact same re-mapping).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11

# Read the current D:G block (columns 4-7) in one shot.
$rng = $ws.Range($ws.Cells.Item(1, 4), $ws.Cells.Item($lastRow, 7))
$vals = $rng.Value2

$newVals = New-Object 'object[,]' $lastRow, 4

for ($r = 1; $r -le $lastRow; $r++) {
    $d = $vals[$r, 1]
    $e = $vals[$r, 2]
    $f = $vals[$r, 3]
    $g = $vals[$r, 4]

    $newVals[$r - 1, 0] = $e   # new D = old category-name
    $newVals[$r - 1, 1] = $g   # new E = old group-name
    $newVals[$r - 1, 2] = $f   # new F = old category-code (unchanged)
    $newVals[$r - 1, 3] = $d   # new G = old group-code
}

# Temporarily force the destination range to Text format so numeric-looking
# codes (e.g. "110", "111") are written back as text, matching the original
# shared-string representation instead of being reinterpreted as numbers.
$rng.NumberFormat = "@"
$rng.Value2 = $newVals

# Restore the plain/default cell style so no visible formatting change (and
# no stray style index) is left behind on the cells themselves.
$rng.Style = "Normal"
